$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 practice date moved from 9/25/2025 to 9/28/2025
$ws.Range("A3").Value = "9/28/2025"

# Row 3 practice time changed from "6:00pm - 8:00pm" to "3:00p - 6:00pm"
$ws.Range("B3").Value = "3:00p - 6:00pm"

# Carry the formatting of the existing data row down into the two new rows
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 4: another practice on 9/30/2025, same time slot as row 2
$ws.Range("A4").Value = "9/30/2025"
$ws.Range("B4").Value = "6:30pm - 8:30pm"

# Row 5 is left blank (freshly formatted, ready for the next entry)

# Grow the table to cover the new rows
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:B5"))

$ws.Range("A5").Select()
